$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, shifting existing rows 18-33 down to 19-34.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the O.R. Tambo International Airport stop.
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 45633
$ws.Range("B18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C18").Value = "공항"
$ws.Range("D18").Value = "O.R. 탐보 국제공항"
$ws.Range("E18").Value = "https://www.google.com/maps/place/O.R.+%ED%83%90%EB%B3%B4+%EA%B5%AD%EC%A0%9C%EA%B3%B5%ED%95%AD/@-26.1393913,28.24422,17z/data=!3m1!4b1!4m6!3m5!1s0x1e95143805a229c3:0xb3bf1c40792821d6!8m2!3d-26.1393913!4d28.2467949!16zL20vMHFuMnY?entry=ttu&g_ep=EgoyMDI0MTEyNC4xIKXMDSoASAFQAw%3D%3D"
$ws.Range("F18").Value = -26.1393913
$ws.Range("G18").Value = 28.24422
